# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F value for sheet "展览"
$exhibitionUpdates = @{
    2  = 1319
    3  = 1198
    4  = 14489
    5  = 17198
    7  = 142
    8  = 51
    9  = 52
    17 = 16
    20 = 1310
    24 = 0
    25 = 7084
    28 = 1156
    29 = 30
    31 = 46
    32 = 5827
    36 = 4995
}

# Row -> new F value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 1319
    3  = 1198
    4  = 14489
    5  = 17198
    7  = 142
    8  = 51
    9  = 52
    17 = 16
    20 = 1310
    26 = 7084
    29 = 1156
    30 = 30
    32 = 46
    34 = 5827
    38 = 4995
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
